$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Issue Log")

$ws.Cells.Item(1, 2).Value = ""

$ws.Range("A2:E5").Value = ""

$ws.Rows("2:5").AutoFit()

$ws.Columns("A").ColumnWidth = 6.5703125
$ws.Columns("C").ColumnWidth = 16.7109375
$ws.Columns("D").ColumnWidth = 13.7109375
$ws.Columns("E").ColumnWidth = 8.42578125

$ws.Range("F16").Select()
